$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update dialogue text cells (column B) with revised script wording
$ws.Range("B5").Value = 'Sir, this young lady’s skills seem quite impressive'
$ws.Range("B6").Value = 'Yao, you’re being too hasty.'
$ws.Range("B7").Value = 'Look at the way she swings the weapon——when it cuts through the air like a feather, it means she’s not applying power correctly.'
$ws.Range("B10").Value = 'I see. You’re amazing——you can spot these details just from a few moves.'
$ws.Range("B13").Value = 'May I ask your name please?'
$ws.Range("B15").Value = 'I’m Chen, the best martial artist in this entire manor.'
$ws.Range("B16").Value = 'My apologies, you must be the top disciple of the manor.'
$ws.Range("B17").Value = 'Top disciple? That’s only because there are just two guards left in Qingliu Manor now.'
$ws.Range("B19").Value = ' <color=#00CC00>(No wonder we’ve walked around so long and seen so few people.)</color>'
$ws.Range("B21").Value = 'Enough chatting——who are you anyway?'
$ws.Range("B22").Value = 'I’m Judge Dee, and this is my student, Yao.'
$ws.Range("B25").Value = 'I really hope I get the chance to spar with them someday.'
$ws.Range("B32").Value = 'Hello!'
$ws.Range("B33").Value = 'Hello——may I ask if you’re a physician?'
$ws.Range("B35").Value = 'I practice medicine in JiuJiang county at the foot of the mountain. I came up a few days ago for a consultation.'
$ws.Range("B36").Value = 'Are you familiar with the Lord?'
$ws.Range("B41").Value = 'Sir, it’s getting late——shall we head back and rest?'

# Clear the now-removed Action2 markers
$ws.Range("J10").ClearContents()
$ws.Range("J40").ClearContents()

# Row heights adjust because of the re-wrapped text in column B
$ws.Rows.Item(19).RowHeight = 34
$ws.Rows.Item(25).RowHeight = 34

# Update the active selection left by the editor
$ws.Range("J10").Select()
